$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 228.57143
$ws.Range("I5").Value = 217
$ws.Range("J5").Value = 237.25
$ws.Range("K5").Value = 217
$ws.Range("L5").Value = 237.25
$ws.Range("M5").Value = -102
$ws.Range("N5").Value = -467.25
$ws.Range("H11").Value = 22.777779
$ws.Range("I11").Value = 22.777779
$ws.Range("K11").Value = 22.777779
$ws.Range("M11").Value = 117.222221
$ws.Range("H17").Value = 475555.6
$ws.Range("I17").Value = 1565.5
$ws.Range("J17").Value = 529725.9
$ws.Range("K17").Value = 4696.5
$ws.Range("L17").Value = 1589177.7
$ws.Range("M17").Value = -4528.5
$ws.Range("N17").Value = -1589513.7
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H33").Value = 5882523
$ws.Range("J33").Value = 272.66666
$ws.Range("L33").Value = 272.66666
$ws.Range("N33").Value = -730.66666
$ws.Range("H38").Value = 104.111115
$ws.Range("I38").Value = 54.625
$ws.Range("K38").Value = 163.875
$ws.Range("M38").Value = 208.125
$ws.Range("H43").Value = 24893.125
$ws.Range("I43").Value = 45500
$ws.Range("J43").Value = 18024.166
$ws.Range("K43").Value = 45500
$ws.Range("L43").Value = 18024.166
$ws.Range("M43").Value = -45431
$ws.Range("N43").Value = -18162.166
$ws.Range("H62").Value = 7653.2104
$ws.Range("I62").Value = 6773.727
$ws.Range("K62").Value = 6773.727
$ws.Range("M62").Value = -6149.727
$ws.Range("H65").Value = 7653.2104
$ws.Range("I65").Value = 6773.727
$ws.Range("K65").Value = 33868.635
$ws.Range("M65").Value = -30748.635
$ws.Range("H80").Value = 1143.1428
$ws.Range("I80").Value = 1329.9286
$ws.Range("J80").Value = 1018.619
$ws.Range("K80").Value = 3989.7858
$ws.Range("L80").Value = 3055.857
$ws.Range("M80").Value = -2991.7858
$ws.Range("N80").Value = -5051.857
$ws.Range("H83").Value = 1143.1428
$ws.Range("I83").Value = 1329.9286
$ws.Range("J83").Value = 1018.619
$ws.Range("K83").Value = 11969.3574
$ws.Range("L83").Value = 9167.571
$ws.Range("M83").Value = -6977.357399999999
$ws.Range("N83").Value = -19151.571
$ws.Range("H86").Value = 4099.7
$ws.Range("I86").Value = 3571
$ws.Range("J86").Value = 5333.3335
$ws.Range("K86").Value = 3571
$ws.Range("L86").Value = 5333.3335
$ws.Range("M86").Value = -2448
$ws.Range("N86").Value = -7579.3335
$ws.Range("H89").Value = 4099.7
$ws.Range("I89").Value = 3571
$ws.Range("J89").Value = 5333.3335
$ws.Range("K89").Value = 17855
$ws.Range("L89").Value = 26666.6675
$ws.Range("M89").Value = -12239
$ws.Range("N89").Value = -37898.6675
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 6264.75
$ws.Range("I107").Value = 3332.5715
$ws.Range("J107").Value = 10369.8
$ws.Range("K107").Value = 3332.5715
$ws.Range("L107").Value = 10369.8
$ws.Range("M107").Value = -1412.5715
$ws.Range("N107").Value = -14209.8
$ws.Range("H116").Value = 14853.444
$ws.Range("I116").Value = 14853.444
$ws.Range("K116").Value = 14853.444
$ws.Range("M116").Value = -11411.444
$ws.Range("H132").Value = 70493.69
$ws.Range("I132").Value = 38637.668
$ws.Range("K132").Value = 115913.004
$ws.Range("M132").Value = -113383.004
$ws.Range("H138").Value = 3962.7913
$ws.Range("I138").Value = 1823.2307
$ws.Range("J138").Value = 4319.385
$ws.Range("K138").Value = 5469.6921
$ws.Range("L138").Value = 12958.155
$ws.Range("M138").Value = -329.6921000000002
$ws.Range("N138").Value = -23238.155
$ws.Range("H141").Value = 1820.25
$ws.Range("I141").Value = 1119
$ws.Range("J141").Value = 3924
$ws.Range("K141").Value = 3357
$ws.Range("L141").Value = 11772
$ws.Range("M141").Value = 1823
$ws.Range("N141").Value = -22132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 2331.1904
$ws.Range("I32").Value = 1927.921
$ws.Range("K32").Value = 1927.921
$ws.Range("M32").Value = -1640.921
$ws.Range("H61").Value = 2403.44
$ws.Range("I61").Value = 1961.1904
$ws.Range("J61").Value = 4725.25
$ws.Range("K61").Value = 1961.1904
$ws.Range("L61").Value = 4725.25
$ws.Range("M61").Value = -1749.1904
$ws.Range("N61").Value = -5149.25
$ws.Range("H74").Value = 7112.5293
$ws.Range("I74").Value = 784.9167
$ws.Range("J74").Value = 22298.8
$ws.Range("K74").Value = 784.9167
$ws.Range("L74").Value = 22298.8
$ws.Range("M74").Value = 89.08330000000001
$ws.Range("N74").Value = -24046.8
$ws.Range("H77").Value = 7112.5293
$ws.Range("I77").Value = 784.9167
$ws.Range("J77").Value = 22298.8
$ws.Range("K77").Value = 3924.5835
$ws.Range("L77").Value = 111494
$ws.Range("M77").Value = 443.4165000000003
$ws.Range("N77").Value = -120230
$ws.Range("H97").Value = 1859.7222
$ws.Range("I97").Value = 1920.8823
$ws.Range("K97").Value = 1920.8823
$ws.Range("M97").Value = -1424.8823
$ws.Range("H102").Value = 4323
$ws.Range("I102").Value = 4252.2
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 4252.2
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -2630.2
$ws.Range("N102").Value = -7744
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 1929.3
$ws.Range("I122").Value = 1765.8889
$ws.Range("K122").Value = 5297.6667
$ws.Range("M122").Value = -2847.6667
$ws.Range("H132").Value = 1726.579
$ws.Range("I132").Value = 1655.8334
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4967.5002
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2437.5002
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2403.44
$ws.Range("I136").Value = 1961.1904
$ws.Range("J136").Value = 4725.25
$ws.Range("K136").Value = 5883.5712
$ws.Range("L136").Value = 14175.75
$ws.Range("M136").Value = -3333.5712
$ws.Range("N136").Value = -19275.75
$ws.Range("H139").Value = 98147.5
$ws.Range("J139").Value = 98147.5
$ws.Range("L139").Value = 98147.5
$ws.Range("N139").Value = -108427.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2131.3684
$ws.Range("I20").Value = 1822.0385
$ws.Range("J20").Value = 2801.5833
$ws.Range("K20").Value = 1822.0385
$ws.Range("L20").Value = 2801.5833
$ws.Range("M20").Value = -1575.0385
$ws.Range("N20").Value = -3295.5833
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H99").Value = 2766.5
$ws.Range("I99").Value = 1092
$ws.Range("J99").Value = 6952.75
$ws.Range("K99").Value = 1092
$ws.Range("L99").Value = 6952.75
$ws.Range("M99").Value = 406
$ws.Range("N99").Value = -9948.75
$ws.Range("H107").Value = 855.05884
$ws.Range("I107").Value = 846.03125
$ws.Range("K107").Value = 846.03125
$ws.Range("M107").Value = 1073.96875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 133.2
$ws.Range("J7").Value = 440
$ws.Range("L7").Value = 440
$ws.Range("N7").Value = -666
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 36106.97
$ws.Range("I31").Value = 45320.26
$ws.Range("K31").Value = 45320.26
$ws.Range("M31").Value = -45025.26
$ws.Range("H34").Value = 36106.97
$ws.Range("I34").Value = 45320.26
$ws.Range("K34").Value = 45320.26
$ws.Range("M34").Value = -45118.26
$ws.Range("H35").Value = 627.2222
$ws.Range("I35").Value = 580.625
$ws.Range("K35").Value = 580.625
$ws.Range("M35").Value = -286.625
$ws.Range("H74").Value = 75000
$ws.Range("J74").Value = 75000
$ws.Range("L74").Value = 75000
$ws.Range("N74").Value = -76748
$ws.Range("H77").Value = 75000
$ws.Range("J77").Value = 75000
$ws.Range("L77").Value = 225000
$ws.Range("N77").Value = -233736
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 714.8333
$ws.Range("J122").Value = 450
$ws.Range("L122").Value = 1350
$ws.Range("N122").Value = -6250
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 20744.414
$ws.Range("I134").Value = 9460.392
$ws.Range("K134").Value = 28381.176
$ws.Range("M134").Value = -25846.176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 25700
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 25700
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 77100
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -77324
$ws.Range("H29").Value = 4208.25
$ws.Range("I29").Value = 900
$ws.Range("J29").Value = 4869.9
$ws.Range("K29").Value = 2700
$ws.Range("L29").Value = 14609.7
$ws.Range("M29").Value = -2423
$ws.Range("N29").Value = -15163.7
$ws.Range("H39").Value = 4884.6875
$ws.Range("I39").Value = 80
$ws.Range("K39").Value = 240
$ws.Range("M39").Value = 54
$ws.Range("H106").Value = 6535.4
$ws.Range("I106").Value = 7599.5
$ws.Range("K106").Value = 22798.5
$ws.Range("M106").Value = -21852.5
$ws.Range("H113").Value = 535.9545000000001
$ws.Range("I113").Value = 478.9
$ws.Range("K113").Value = 1436.7
$ws.Range("M113").Value = 733.3000000000002
$ws.Range("H135").Value = 25700
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 25700
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 231300
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -236370
$ws.Range("H140").Value = 2396.9375
$ws.Range("I140").Value = 2396.9375
$ws.Range("K140").Value = 7190.8125
$ws.Range("M140").Value = -2010.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 10043
$ws.Range("J53").Value = 10043
$ws.Range("L53").Value = 10043
$ws.Range("N53").Value = -11305
$ws.Range("H70").Value = 5893.8887
$ws.Range("I70").Value = 4843.4346
$ws.Range("K70").Value = 4843.4346
$ws.Range("M70").Value = -4573.4346
$ws.Range("H73").Value = 5893.8887
$ws.Range("I73").Value = 4843.4346
$ws.Range("K73").Value = 4843.4346
$ws.Range("M73").Value = -3907.4346
$ws.Range("H80").Value = 2248.5625
$ws.Range("I80").Value = 1743.8
$ws.Range("K80").Value = 1743.8
$ws.Range("M80").Value = -745.8
$ws.Range("H83").Value = 2248.5625
$ws.Range("I83").Value = 1743.8
$ws.Range("K83").Value = 8719
$ws.Range("M83").Value = -3727
$ws.Range("H93").Value = 26250
$ws.Range("J93").Value = 26250
$ws.Range("L93").Value = 26250
$ws.Range("N93").Value = -29994
$ws.Range("H97").Value = 946.25
$ws.Range("I97").Value = 821.25
$ws.Range("K97").Value = 821.25
$ws.Range("M97").Value = -325.25
$ws.Range("H126").Value = 17780.408
$ws.Range("J126").Value = 4630.5
$ws.Range("L126").Value = 13891.5
$ws.Range("N126").Value = -18831.5
$ws.Range("H132").Value = 458068.47
$ws.Range("I132").Value = 458068.47
$ws.Range("K132").Value = 1374205.41
$ws.Range("M132").Value = -1371675.41

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3745.9773
$ws.Range("J16").Value = 6070.5454
$ws.Range("L16").Value = 6070.5454
$ws.Range("N16").Value = -6410.5454
$ws.Range("H22").Value = 1132.75
$ws.Range("I22").Value = 1272.4
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1272.4
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -977.4000000000001
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 1132.75
$ws.Range("I27").Value = 1272.4
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 1272.4
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -1165.4
$ws.Range("N27").Value = -1114
$ws.Range("H40").Value = 6084.3335
$ws.Range("I40").Value = 5701.2
$ws.Range("K40").Value = 5701.2
$ws.Range("M40").Value = -5565.2
$ws.Range("H46").Value = 2670
$ws.Range("I46").Value = 499
$ws.Range("K46").Value = 499
$ws.Range("M46").Value = -311
$ws.Range("H55").Value = 174.88889
$ws.Range("I55").Value = 173.6
$ws.Range("J55").Value = 181.33333
$ws.Range("K55").Value = 173.6
$ws.Range("L55").Value = 181.33333
$ws.Range("M55").Value = -0.5999999999999943
$ws.Range("N55").Value = -527.3333299999999
$ws.Range("H61").Value = 3205.9285
$ws.Range("I61").Value = 2953.182
$ws.Range("K61").Value = 2953.182
$ws.Range("M61").Value = -2751.182
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H75").Value = 41000
$ws.Range("J75").Value = 41000
$ws.Range("L75").Value = 41000
$ws.Range("N75").Value = -42872
$ws.Range("H78").Value = 41000
$ws.Range("J78").Value = 41000
$ws.Range("L78").Value = 123000
$ws.Range("N78").Value = -132360
$ws.Range("H100").Value = 59391.367
$ws.Range("I100").Value = 67764.75
$ws.Range("K100").Value = 67764.75
$ws.Range("M100").Value = -67223.75
$ws.Range("H113").Value = 3205.9285
$ws.Range("I113").Value = 2953.182
$ws.Range("K113").Value = 2953.182
$ws.Range("M113").Value = -783.1819999999998
$ws.Range("H122").Value = 7541.737
$ws.Range("I122").Value = 5956.7144
$ws.Range("J122").Value = 11979.8
$ws.Range("K122").Value = 17870.1432
$ws.Range("L122").Value = 35939.39999999999
$ws.Range("M122").Value = -15420.1432
$ws.Range("N122").Value = -40839.39999999999
$ws.Range("H132").Value = 6410.3335
$ws.Range("I132").Value = 6223.5
$ws.Range("J132").Value = 6559.8
$ws.Range("K132").Value = 18670.5
$ws.Range("L132").Value = 19679.4
$ws.Range("M132").Value = -16140.5
$ws.Range("N132").Value = -24739.4
$ws.Range("H136").Value = 4164.5
$ws.Range("I136").Value = 4260.2354
$ws.Range("J136").Value = 3839
$ws.Range("K136").Value = 12780.7062
$ws.Range("L136").Value = 11517
$ws.Range("M136").Value = -10230.7062
$ws.Range("N136").Value = -16617

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 33100.2
$ws.Range("I2").Value = 25249.5
$ws.Range("J2").Value = 38334
$ws.Range("K2").Value = 25249.5
$ws.Range("L2").Value = 38334
$ws.Range("M2").Value = -25137.5
$ws.Range("N2").Value = -38558
$ws.Range("H81").Value = 9081.9
$ws.Range("I81").Value = 16376.857
$ws.Range("J81").Value = 5153.846
$ws.Range("K81").Value = 32753.714
$ws.Range("L81").Value = 10307.692
$ws.Range("M81").Value = -31692.714
$ws.Range("N81").Value = -12429.692
$ws.Range("H84").Value = 9081.9
$ws.Range("I84").Value = 16376.857
$ws.Range("J84").Value = 5153.846
$ws.Range("K84").Value = 163768.57
$ws.Range("L84").Value = 51538.45999999999
$ws.Range("M84").Value = -158464.57
$ws.Range("N84").Value = -62146.45999999999
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 1143.6389
$ws.Range("I100").Value = 968.3182
$ws.Range("J100").Value = 1419.1428
$ws.Range("K100").Value = 1936.6364
$ws.Range("L100").Value = 2838.2856
$ws.Range("M100").Value = -1395.6364
$ws.Range("N100").Value = -3920.2856
$ws.Range("H122").Value = 3421.4546
$ws.Range("I122").Value = 3085.8333
$ws.Range("K122").Value = 9257.499899999999
$ws.Range("M122").Value = -6807.499899999999
$ws.Range("H126").Value = 5337.2
$ws.Range("I126").Value = 4921.5
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 14764.5
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -12294.5
$ws.Range("N126").Value = -25940
$ws.Range("H136").Value = 2402.1035
$ws.Range("I136").Value = 2165.1667
$ws.Range("J136").Value = 3539.4
$ws.Range("K136").Value = 6495.500100000001
$ws.Range("L136").Value = 10618.2
$ws.Range("M136").Value = -3945.500100000001
$ws.Range("N136").Value = -15718.2
